# Generate Report for Handback
#
# Row 7 in both the "zh-cn" and "de-de" sheets corresponds to the
# 08814469-19cd-4a3e-a654-ed96e2a330a5.md file, which was previously
# "Ready for handoff" (no handback processed yet: Latest Target File,
# Latest Handback File, Latest Handback DateTime and Error Detail were
# blank / placeholder). This change fills in the handback report:
#   - Latest Target File (I7): link to the source .md file
#   - Latest Handback File (J7): the generated .xlf handback file name
#   - Latest Handback DateTime (K7): the timestamp the handback was processed
#   - Error Detail (P7): the file was stale, so an error message is recorded

$wb = $excel.ActiveWorkbook

$mdTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5634e97dbcfabab657511e9e93046866e2add140/e2e/08814469-19cd-4a3e-a654-ed96e2a330a5.md"
$mdDisplay = "08814469-19cd-4a3e-a654-ed96e2a330a5.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cde79be54fd64e8eb84d7be62b9b4cf5e02d55a2/e2e/08814469-19cd-4a3e-a654-ed96e2a330a5.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5634e97dbcfabab657511e9e93046866e2add140/e2e/08814469-19cd-4a3e-a654-ed96e2a330a5.md."

# ---------------- zh-cn sheet ----------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("J7").Value = "08814469-19cd-4a3e-a654-ed96e2a330a5.8614b93c49bdfd985efe94faedb86cfd843a0f6b.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-16 10:55:46"
$wsZh.Range("P7").Value = $errorDetail

# Latest Target File gets a hyperlink back to the source markdown file
# (same as column A), which also sets the cell's display text/value and
# applies the hyperlink style.
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $mdTarget, "", "", $mdDisplay)

# ---------------- de-de sheet ----------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("J7").Value = "08814469-19cd-4a3e-a654-ed96e2a330a5.8614b93c49bdfd985efe94faedb86cfd843a0f6b.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-16 10:55:54"
$wsDe.Range("P7").Value = $errorDetail

$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $mdTarget, "", "", $mdDisplay)
